# "Settings back to qfin22 presentation - recalibration"
# Updates the recalibrated parameter values on the "Linear" and "NonLinear"
# worksheets (including the abs_epsi_autocorr array strings) to reflect the
# new calibration run.

$wb = $excel.ActiveWorkbook

# --- Linear sheet ---
$wsLinear = $wb.Worksheets.Item("Linear")
$wsLinear.Range("B2").Value = -0.1234521727377149
$wsLinear.Range("B3").Value = 0.1394397578864278
$wsLinear.Range("B4").Value = 879.0626802071353
$wsLinear.Range("B5").Value = "[0.9999999999999999, 0.21262385955974378, 0.07529032599515867, 0.10340021105703882, 0.06688382568203936, 0.04187258675481803, 0.19296588660627212, 0.32560474823588526, 0.17572974856412568, 0.03528694163264741, 0.05406449694684162, 0.0577469080820669, 0.012877242238997404, 0.15779499682617557, 0.3100053258955958, 0.15144080541633376, 0.0010807386165646916, 0.04017423819544801, 0.026467968265428847, -0.007158722646460934]"

# --- NonLinear sheet ---
$wsNonLinear = $wb.Worksheets.Item("NonLinear")
$wsNonLinear.Range("B4").Value = -0.4306139585734692
$wsNonLinear.Range("B5").Value = 0.1535499012469229
$wsNonLinear.Range("B6").Value = 918.8977458070067
$wsNonLinear.Range("B7").Value = 1.306372138250548
$wsNonLinear.Range("B8").Value = 0.02828024968310232
$wsNonLinear.Range("B9").Value = 843.0544412962871
$wsNonLinear.Range("B10").Value = "[0.9999999999999999, 0.2092967758037653, 0.07249458564993277, 0.10480220818339674, 0.06819853810308928, 0.040672243371456716, 0.19127616741882666, 0.32181277231045924, 0.17389226724091802, 0.0355946698281903, 0.05521503062263752, 0.05897915828522973, 0.013413487042602073, 0.15580651845556104, 0.3072501473127087, 0.14990690995558545, 0.0017331863438872237, 0.040969868963442915, 0.026661197257683133, -0.007116626777000597]"
